# The edit re-orders the 19 data records (rows 2-20) of the "Artfynd"
# sheet into a new row order; every cell in a source row (A:AY) moves as
# a unit to its destination row, values/types unchanged, header row (1)
# untouched.
#
# before-row -> after-row mapping (derived from the Id column, A, which
# travels with its record):
#   2->17  3->9   4->11  5->12  6->4   7->5   8->19  9->18  10->16
#   11->13 12->14 13->3  14->10 15->8  16->7  17->15 18->2  19->20 20->6
#
# Because this permutation contains cycles, a direct row-to-row copy
# would clobber source data before it gets a chance to move elsewhere.
# So: stage every source row far below the used range first, clear the
# originals, then copy each staged row into its real destination, and
# finally wipe the staging area.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    2 = 17
    3 = 9
    4 = 11
    5 = 12
    6 = 4
    7 = 5
    8 = 19
    9 = 18
    10 = 16
    11 = 13
    12 = 14
    13 = 3
    14 = 10
    15 = 8
    16 = 7
    17 = 15
    18 = 2
    19 = 20
    20 = 6
}

$firstRow = 2
$lastRow = 20
$lastCol = "AY"
$stageOffset = 1000

# 1) Stage every source row (its full A:AY extent) well below the data.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $src = $ws.Range("A" + $r + ":" + $lastCol + $r)
    $stage = $ws.Range("A" + ($r + $stageOffset) + ":" + $lastCol + ($r + $stageOffset))
    $src.Copy($stage)
}

# 2) Clear the original rows completely so stale values can't survive in
#    spots whose incoming row has blank cells there (Copy does not blank
#    out a destination cell when the source cell is empty).
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Range("A" + $r + ":" + $lastCol + $r).ClearContents()
}

# 3) Move staged rows into their real destinations.
foreach ($r in $map.Keys) {
    $dst = $map[$r]
    $stage = $ws.Range("A" + ($r + $stageOffset) + ":" + $lastCol + ($r + $stageOffset))
    $target = $ws.Range("A" + $dst + ":" + $lastCol + $dst)
    $stage.Copy($target)
}

# 4) Clean up the staging area.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Range("A" + ($r + $stageOffset) + ":" + $lastCol + ($r + $stageOffset)).ClearContents()
}
